# Weekly update: insert a new week's price report as row 79 on the
# "Haba" (Feria Lagunitas de Puerto Montt) sheet, pushing the existing
# rows 79:126 down to 80:127.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 79, shifting rows 79:126 -> 80:127.
$ws.Rows.Item(79).Insert()

# Populate the new row 79 with this week's data (same market/region as
# the former row 79, with updated date, volume and prices).
$ws.Range("A79").Value = 4
$ws.Range("B79").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C79").Value = "Los Lagos"
$ws.Range("D79").Value = 45062
$ws.Range("E79").Value = 10
$ws.Range("F79").Value = 100112026
$ws.Range("G79").Value = "Haba"
$ws.Range("H79").Value = "Sin especificar"
$ws.Range("I79").Value = "Primera"
$ws.Range("J79").Value = 120
$ws.Range("K79").Value = 21000
$ws.Range("L79").Value = 21000
$ws.Range("M79").Value = 21000
$ws.Range("N79").Value = "$/saco 25 kilos"
$ws.Range("O79").Value = "Provincia de Limarí"
$ws.Range("P79").Value = 840
$ws.Range("Q79").Value = 25
$ws.Range("R79").Value = "Hortaliza"
